# Blog export sheet refresh ("fix giao dien blog")
# Replaces the old sample blog rows with the current blog list and
# removes the now-unused trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Helper: write a plain-text value into a cell without letting the
# "smart" input layer reinterpret date-shaped strings (e.g.
# "2022-12-08") as a date serial number / date-formatted cell. We
# build the text as a formula result, then freeze it back into a
# literal value in place via copy / paste-special(values) - this
# keeps the destination cell's existing style untouched.
# ------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163) # xlPasteValues
}

# Drop the old rows 7-9 - the refreshed data only spans down to row 6.
$ws.Rows("7:9").Delete()

# Row 2
$ws.Range("A2").Value = 8
$ws.Range("B2").Value = "Bóng đá"
$ws.Range("C2").Value = "NÔI DUNG BLOG"
$ws.Range("D2").Value = 1
Set-TextValue $ws.Range("E2") "2022-12-08"
$ws.Range("F2").Value = "75258550_2545915118988163_3298790051630022656_n.png"

# Row 3
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "Bóng Chuyền"
$ws.Range("C3").Value = "NÔI DUNG BLOG"
$ws.Range("D3").Value = 2
Set-TextValue $ws.Range("E3") "2022-12-08"
$ws.Range("F3").Value = "311038288_555308266596129_7246856222758665646_n.jpeg"

# Row 4
$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Game"
$ws.Range("C4").Value = "NÔI DUNG BLOG"
$ws.Range("D4").Value = 2
Set-TextValue $ws.Range("E4") "2022-12-08"
$ws.Range("F4").Value = "20*7.png"

# Row 5
$ws.Range("A5").Value = 12
$ws.Range("B5").Value = "Cầu Long"
$ws.Range("C5").Value = "Nội dung"
$ws.Range("D5").ClearContents()
Set-TextValue $ws.Range("E5") "2022-12-08"
$ws.Range("F5").Value = "Screenshot 2022-10-31 at 15.44.44.png"

# Row 6
$ws.Range("A6").Value = 13
$ws.Range("B6").Value = "Cầu Long"
$ws.Range("C6").Value = "NÔI DUNG BLOG"
$ws.Range("D6").Value = 2
Set-TextValue $ws.Range("E6") "2022-12-08"
$ws.Range("F6").Value = "Dương Tùng.png"

$ws.Range("F6").Select()
